# Append a new row of Kaspa buy data (run on 2025-10-03) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-like string. Force a text number format first so
# Excel stores the literal string "10/03/2025" instead of auto-converting
# it to a date serial number, then restore the default "Normal" style so
# no extra formatting is left behind on the cell.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "10/03/2025"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = 304.5709999999999
$ws.Range("C8").Value = 0.1641653341913709
$ws.Range("D8").Value = 25
